$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 474, shifting existing rows 474:570 down to 475:571.
$ws.Rows.Item(474).Insert()

# Populate the newly inserted row 474 with the new record.
$ws.Cells.Item(474, 1).Value = 9
$ws.Cells.Item(474, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(474, 3).Value = "Metropolitana"
$ws.Cells.Item(474, 4).Value = 44785
$ws.Cells.Item(474, 5).Value = 13
$ws.Cells.Item(474, 6).Value = 100112031
$ws.Cells.Item(474, 7).Value = "Poroto verde"
$ws.Cells.Item(474, 8).Value = "Magnum"
$ws.Cells.Item(474, 9).Value = "Primera"
$ws.Cells.Item(474, 10).Value = 52
$ws.Cells.Item(474, 11).Value = 38000
$ws.Cells.Item(474, 12).Value = 38000
$ws.Cells.Item(474, 13).Value = 38000
$ws.Cells.Item(474, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(474, 15).Value = "Perú"
$ws.Cells.Item(474, 16).Value = 1520
$ws.Cells.Item(474, 17).Value = 25
$ws.Cells.Item(474, 18).Value = "Hortaliza"
